$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the stable-coin address text in H2 ---
$ws.Range("H2").Value = "0x8ac76a51cc950d9822d68b83fe1ad97b32cd580d"

# Clear H2's highlighted/bordered header look: no fill, no border,
# general horizontal alignment, wrap text enabled.
$ws.Range("H2").Interior.Pattern = -4142
$ws.Range("H2").Borders.LineStyle = -4142
$ws.Range("H2").HorizontalAlignment = 1
$ws.Range("H2").WrapText = $true

# --- "combien de dollars j'ajoute a chaque achat" value ---
$ws.Range("D2").Value = 2

# --- First data row (row 4): new price / amount / date ---
$ws.Range("A4").Value = 0.2675
$ws.Range("B4").Value = 10
$purchaseDate = Get-Date -Year 2025 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("C4").Value = $purchaseDate

# Remove the old lime-green highlight fill from A4:C4 (keep their other
# formatting - wrap/center/date format - untouched).
$ws.Range("A4:C4").Interior.Pattern = -4142
$ws.Range("A4:C4").WrapText = $true
$ws.Range("A4:C4").HorizontalAlignment = -4108

# D4 instead now receives the green highlight (same tone as the row-2
# header band).
$ws.Range("D4").Interior.Color = 0x6BEAD4

# --- Update the selection shown when the sheet is opened ---
$ws.Range("A5:D22").Select()
$ws.Range("D22").Activate()
